# Add a "Credit" note after the DFS reflection paragraph, crediting ChatGPT
# for the idea of using a while loop instead of recursion (to avoid hitting
# the recursion stack limit).

$d = $word.ActiveDocument

# Locate the end of the sentence the new content should be appended after.
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "testing and documentation in software development.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the anchor sentence to attach the credit note to."
}

# Collapse to the end of the found text so we insert right after the period,
# before the following (bold, larger) line-break run.
$anchor.Collapse(0)

# A manual line break (Chr(11) / vertical-tab) is how Word represents <w:br/>
# when it's embedded inside a run's text. Inserting it together with the new
# sentence keeps it in the same run/formatting context (matching the size
# 11pt / non-bold styling of the paragraph it is being appended to).
$lineBreak = [char]11
$creditText = "Credit: Used ChatGPT to get the alternative idea of using the while loop as I was hitting the recursion stack Limit."

$anchor.InsertAfter($lineBreak + $creditText)
